$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Alone.menage"
$ws.Range("C1").Value = "With the family of origin (parents, etc.).menage"
$ws.Range("D1").Value = "With partner / children.menage"
$ws.Range("E1").Value = "With friends or other people (with no family relation).menage"
$ws.Range("F1").Value = "In detention.menage"
$ws.Range("G1").Value = "In institutions /shelters (not detention).menage"
$ws.Range("H1").Value = "Other.menage"
$ws.Range("I1").Value = "Not known / missing.menage"
$ws.Range("J1").Value = "Total.menage"
